# Refresh the cryptos price/volume snapshot (cols D=Price, E=Volume(1h)).
# Rows 45/46 additionally swap coin identity (B=Coin, C=Link) because the
# ranking order of FTXToken and InjectiveProtocol flipped in this update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.231.19'
$ws.Cells.Item(2, 5).Value = '  +0.21%  '
$ws.Cells.Item(3, 4).Value = '2.023.81'
$ws.Cells.Item(3, 5).Value = '  -0.01%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 4).Value = '228.59'
$ws.Cells.Item(5, 5).Value = '  +0.62%  '
$ws.Cells.Item(6, 5).Value = '  +0.46%  '
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 4).Value = '56.09'
$ws.Cells.Item(8, 5).Value = '  +1.31%  '
$ws.Cells.Item(9, 4).Value = '0.378'
$ws.Cells.Item(9, 5).Value = '  -1.12%  '
$ws.Cells.Item(10, 4).Value = '0.0782'
$ws.Cells.Item(10, 5).Value = '  -1.33%  '
$ws.Cells.Item(11, 5).Value = '  -2.58%  '
$ws.Cells.Item(12, 4).Value = '2.321.77'
$ws.Cells.Item(12, 5).Value = '  +0.05%  '
$ws.Cells.Item(13, 4).Value = '14.29'
$ws.Cells.Item(13, 5).Value = '  -0.37%  '
$ws.Cells.Item(14, 4).Value = '20.07'
$ws.Cells.Item(14, 5).Value = '  -2.15%  '
$ws.Cells.Item(15, 5).Value = '  +0.60%  '
$ws.Cells.Item(16, 4).Value = '0.739'
$ws.Cells.Item(16, 5).Value = '  -0.80%  '
$ws.Cells.Item(17, 4).Value = '2.023.53'
$ws.Cells.Item(17, 5).Value = '  +0.05%  '
$ws.Cells.Item(18, 4).Value = '37.153.62'
$ws.Cells.Item(18, 5).Value = '  +0.31%  '
$ws.Cells.Item(19, 4).Value = '6.16'
$ws.Cells.Item(19, 5).Value = '  +1.58%  '
$ws.Cells.Item(20, 4).Value = '69.01'
$ws.Cells.Item(20, 5).Value = '  +0.29%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0818'
$ws.Cells.Item(21, 5).Value = '  -2.38%  '
$ws.Cells.Item(22, 4).Value = '223.32'
$ws.Cells.Item(22, 5).Value = '  +0.01%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  -0.14%  '
$ws.Cells.Item(24, 5).Value = '  +1.68%  '
$ws.Cells.Item(25, 4).Value = '2.21'
$ws.Cells.Item(25, 5).Value = '  -2.29%  '
$ws.Cells.Item(26, 4).Value = '163.37'
$ws.Cells.Item(26, 5).Value = '  -2.39%  '
$ws.Cells.Item(27, 5).Value = '  -3.59%  '
$ws.Cells.Item(28, 5).Value = '  +1.79%  '
$ws.Cells.Item(29, 4).Value = '18.74'
$ws.Cells.Item(29, 5).Value = '  -0.06%  '
$ws.Cells.Item(30, 5).Value = '  -1.32%  '
$ws.Cells.Item(31, 4).Value = '0.118'
$ws.Cells.Item(31, 5).Value = '  +0.33%  '
$ws.Cells.Item(32, 4).Value = '4.46'
$ws.Cells.Item(32, 5).Value = '  -0.61%  '
$ws.Cells.Item(33, 4).Value = '0.0601'
$ws.Cells.Item(33, 5).Value = '  -1.14%  '
$ws.Cells.Item(34, 4).Value = '4.46'
$ws.Cells.Item(34, 5).Value = '  -0.09%  '
$ws.Cells.Item(35, 4).Value = '2.36'
$ws.Cells.Item(35, 5).Value = '  -0.42%  '
$ws.Cells.Item(36, 5).Value = '  +2.95%  '
$ws.Cells.Item(37, 5).Value = '  +0.04%  '
$ws.Cells.Item(38, 4).Value = '3.19'
$ws.Cells.Item(38, 5).Value = '  +0.32%  '
$ws.Cells.Item(39, 4).Value = '5.49'
$ws.Cells.Item(39, 5).Value = '  +1.93%  '
$ws.Cells.Item(40, 4).Value = '1.473.83'
$ws.Cells.Item(40, 5).Value = '  -2.08%  '
$ws.Cells.Item(41, 5).Value = '  -2.64%  '
$ws.Cells.Item(42, 4).Value = '94.66'
$ws.Cells.Item(42, 5).Value = '  -0.60%  '
$ws.Cells.Item(43, 4).Value = '2.81'
$ws.Cells.Item(43, 5).Value = '  -1.48%  '
$ws.Cells.Item(44, 4).Value = '0.0914'
$ws.Cells.Item(44, 5).Value = '  -1.56%  '
$ws.Cells.Item(45, 2).Value = 'FTXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(45, 4).Value = '4.21'
$ws.Cells.Item(45, 5).Value = '  +15.70%  '
$ws.Cells.Item(46, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(46, 4).Value = '16.24'
$ws.Cells.Item(46, 5).Value = '  -1.87%  '
$ws.Cells.Item(47, 5).Value = '  -1.37%  '
$ws.Cells.Item(48, 5).Value = '  +0.19%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '7.10'
$ws.Cells.Item(49, 5).Value = '  -0.81%  '
$ws.Cells.Item(50, 5).Value = '  +0.76%  '
$ws.Cells.Item(51, 4).Value = '2.207.84'
$ws.Cells.Item(51, 5).Value = '  -0.08%  '
